$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the "MuSCs" sending-cluster row) entirely; this shifts the
# old row 4 ("Resolving-Mac") up to become the new row 3, and Excel will
# drop the now-unused "MuSCs" shared string automatically.
$ws.Rows("3").Delete()

# Update the recalculated TPM-derived values for row 2 (ECs -> Ccl3/Ackr2/FAPs)
$ws.Range("G2").Value = 0.4347226666666668
$ws.Range("H2").Value = 1.304168
$ws.Range("I2").Value = 0.00308384331793249
$ws.Range("J2").Value = 0.00308384331793249
$ws.Range("Q2").Value = 0.1193867266862222
$ws.Range("R2").Value = 1.074480540176
$ws.Range("S2").Value = 0.00308384331793249
$ws.Range("T2").Value = 0.00308384331793249

# Update the recalculated TPM-derived values for row 3 (Resolving-Mac -> Ccl3/Ackr2/FAPs)
$ws.Range("G3").Value = 140.5330963333333
$ws.Range("H3").Value = 421.599289
$ws.Range("I3").Value = 0.9969161566820676
$ws.Range("J3").Value = 0.9969161566820675
$ws.Range("Q3").Value = 38.59422949109977
$ws.Range("R3").Value = 347.348065419898
$ws.Range("S3").Value = 0.9969161566820676
$ws.Range("T3").Value = 0.9969161566820675
